# Cmaes_210E_program.xlsx — "working great with beats as pems, interface is
# done except some probable bugs"
#
# Updates a handful of the CMAES run-config values on the "cmaes" sheet and
# moves the workbook's on-screen selection/scroll position to reflect where
# the author was last working (row 34, insigma).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Config value edits -----------------------------------------------
# B14 temp.erfStruct.performance_calculators:
#   struct('CongestionPattern',0.9999,'TVH',0.0001)
#   -> struct('CongestionPattern',[0.9989,1],'TVH',[0.001,2],'TVM',[0.0001,2])
$ws.Range("B14").Value = "struct('CongestionPattern',[0.9989,1],'TVH',[0.001,2],'TVM',[0.0001,2])"

# B22 knobs.underevaluation_tolerance_coefficient: [0.5] -> [0.9]
$ws.Range("B22").Value = "[0.9]"

# B23 knobs.overevaluation_tolerance_coefficient: [1.5] -> [1.1]
$ws.Range("B23").Value = "[1.1]"

# B26 maxEval: [6000] -> [100]
$ws.Range("B26").Value = "[100]"

# B27 maxIter: [6000] -> [100]
$ws.Range("B27").Value = "[100]"

# B34 insigma: [2] -> [1.5]
$ws.Range("B34").Value = "[1.5]"

# --- View / selection state --------------------------------------------
# Author ended up with the window scrolled down to row ~10 and the
# selection on B34 (previously scrolled to row ~7 with selection on B28).
$win = $excel.ActiveWindow
[void]$ws.Range("B34").Select()
$win.ScrollRow = 10
$win.ScrollColumn = 1
